$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Total" header in column X (row 1)
$ws.Range("X1").Value = "Total"

# Add a "Total" column (X) with the row sum for each existing data row (2-6)
$ws.Range("X2").Value = 1978
$ws.Range("X3").Value = 274
$ws.Range("X4").Value = 1086
$ws.Range("X5").Value = 249
$ws.Range("X6").Value = 1333

# New row 7: "Outros" category
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 134
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 40
$ws.Range("F7").Value = 64
$ws.Range("G7").Value = 73
$ws.Range("H7").Value = 83
$ws.Range("I7").Value = 118
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 88
$ws.Range("L7").Value = 105
$ws.Range("M7").Value = 118
$ws.Range("N7").Value = 78
$ws.Range("O7").Value = 103
$ws.Range("P7").Value = 126
$ws.Range("Q7").Value = 128
$ws.Range("R7").Value = 160
$ws.Range("S7").Value = 183
$ws.Range("T7").Value = 106
$ws.Range("U7").Value = 35
$ws.Range("V7").Value = 12
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = 1839

# New row 8: "Total" category (sum over all category rows)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 148
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 14
$ws.Range("E8").Value = 54
$ws.Range("F8").Value = 76
$ws.Range("G8").Value = 95
$ws.Range("H8").Value = 118
$ws.Range("I8").Value = 178
$ws.Range("J8").Value = 158
$ws.Range("K8").Value = 267
$ws.Range("L8").Value = 352
$ws.Range("M8").Value = 488
$ws.Range("N8").Value = 471
$ws.Range("O8").Value = 598
$ws.Range("P8").Value = 637
$ws.Range("Q8").Value = 790
$ws.Range("R8").Value = 827
$ws.Range("S8").Value = 785
$ws.Range("T8").Value = 489
$ws.Range("U8").Value = 168
$ws.Range("V8").Value = 35
$ws.Range("W8").Value = 1
$ws.Range("X8").Value = 6759
